$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")
$ws.Range("D2:D10").Copy()
$ws.Range("E2:E10").PasteSpecial(-4122)
